$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# The keyword list in column A ran from row 1 to row 100. This manual test
# trims the last batch of 14 keywords (rows 87-100), clearing their
# contents while leaving the cell styling (border) in place.
$ws.Range("A87:A100").ClearContents()

# Leave the freshly-cleared block selected and scroll the window down so
# it is in view, matching the state the sheet was left in after the edit.
$ws.Range("A87:A100").Select()
$excel.ActiveWindow.ScrollRow = 82
